# ---------------------------------------------------------------------------
# Applies the two content changes described by the commit's OOXML diff:
#
#  1. Slide 5's table switches its Table Style (a:tblPr/a:tableStyleId) from
#     {87763A4E-454D-410C-A30E-9F54995857A7} to
#     {45E8FC6C-C3A5-4660-97E2-8F6652C1DC9E}.
#
#  2. The presentation's applied theme ("Integral" / Red Violet) reverts to
#     the stock "Office Theme" colour palette. (Font scheme and format
#     scheme -- gradients/lines/effects -- are already identical between the
#     two themes in this deck, so only the 12 theme colours actually need to
#     change.)
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{45E8FC6C-C3A5-4660-97E2-8F6652C1DC9E}", $true)
    }
}

# --- 2. Switch the deck's theme colours back to the default Office palette -
$tcs = $p.Slides.Item(1).ThemeColorScheme

function Set-ThemeColor {
    param($Index, $R, $G, $B)
    $color = $tcs.Colors($Index)
    $color.RGB = $R + ($G * 256) + ($B * 65536)
}

Set-ThemeColor 1  0x00 0x00 0x00   # dk1
Set-ThemeColor 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor 11 0x05 0x63 0xC1   # hyperlink
Set-ThemeColor 12 0x95 0x4F 0x72   # followed hyperlink
